# "summary (windows) max post adj.xlsx" update:
# A second batch of simulation runs finished, filling in the previously
# empty Avg.Cycle.Usage / Survival.Rate / Correctly.Scheduled /
# Correctly.Unscheduled columns (F:I) for the AR1_state_based_logistic /
# Window.Size=36 / Sample.Size=3000 rows that were still pending.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F101").Value = 0.094935300383998764
$ws.Range("G101").Value = 0.93663696968958332
$ws.Range("H101").Value = 0.89961247462631477
$ws.Range("I101").Value = 0.15862077774442621
$ws.Range("F105").Value = 0.10574475609049247
$ws.Range("G105").Value = 0.93663696968958332
$ws.Range("H105").Value = 0.89988958410011044
$ws.Range("I105").Value = 0.15863597460495893
$ws.Range("F109").Value = 0.1231944614553447
$ws.Range("G109").Value = 0.93637283712763231
$ws.Range("H109").Value = 0.90093998553868404
$ws.Range("I109").Value = 0.15870976813150292
$ws.Range("F113").Value = 0.18837019872599728
$ws.Range("G113").Value = 0.93608141748425833
$ws.Range("H113").Value = 0.90687100033681378
$ws.Range("I113").Value = 0.15909061783137471
$ws.Range("F117").Value = 0.24691685334691088
$ws.Range("G117").Value = 0.97677428092060214
$ws.Range("H117").Value = 0.87509875177753205
$ws.Range("I117").Value = 0.15936575109480253
$ws.Range("F121").Value = 0.25685844046940376
$ws.Range("G121").Value = 0.9767590586252558
$ws.Range("H121").Value = 0.87540262392960955
$ws.Range("I121").Value = 0.15940805648786435
$ws.Range("F125").Value = 0.27149235265058214
$ws.Range("G125").Value = 0.97643508169823956
$ws.Range("H125").Value = 0.87601435705368291
$ws.Range("I125").Value = 0.15948090948090948
$ws.Range("F129").Value = 0.33979117670775538
$ws.Range("G129").Value = 0.97571871986499514
$ws.Range("H129").Value = 0.87917197929948254
$ws.Range("I129").Value = 0.15988753388445318
$ws.Range("F133").Value = 0.336597177186369
$ws.Range("G133").Value = 0.98239574368202798
$ws.Range("H133").Value = 0.86721302080348439
$ws.Range("I133").Value = 0.15955614281464672
$ws.Range("F137").Value = 0.34683351820578401
$ws.Range("G137").Value = 0.98236954835007917
$ws.Range("H137").Value = 0.86746299365605528
$ws.Range("I137").Value = 0.15959530251967088
$ws.Range("F141").Value = 0.36189013892909327
$ws.Range("G141").Value = 0.98220329587917288
$ws.Range("H141").Value = 0.8676303209651719
$ws.Range("I141").Value = 0.15962971169025397
$ws.Range("F145").Value = 0.43362715650175115
$ws.Range("G145").Value = 0.9812908406065789
$ws.Range("H145").Value = 0.86957735469823016
$ws.Range("I145").Value = 0.15993587296740208
$ws.Range("F149").Value = 0.52223206158863145
$ws.Range("G149").Value = 0.95361465760238762
$ws.Range("H149").Value = 0.84159954107608093
$ws.Range("I149").Value = 0.15607133643467624
$ws.Range("F153").Value = 0.52980843169437308
$ws.Range("G153").Value = 0.95352420168299457
$ws.Range("H153").Value = 0.841744815901878
$ws.Range("I153").Value = 0.15612067314497929
$ws.Range("F157").Value = 0.54087895307085443
$ws.Range("G157").Value = 0.95324680709902143
$ws.Range("H157").Value = 0.84171632896305126
$ws.Range("I157").Value = 0.15610995418575593
$ws.Range("F161").Value = 0.59383957260441089
$ws.Range("G161").Value = 0.95223535759495859
$ws.Range("H161").Value = 0.84207536451429654
$ws.Range("I161").Value = 0.15623121116975081
$ws.Range("F165").Value = 0.60695249761999315
$ws.Range("G165").Value = 0.93928945437447908
$ws.Range("H165").Value = 0.84801903907301568
$ws.Range("I165").Value = 0.15884483054444315
$ws.Range("F169").Value = 0.61259306516572287
$ws.Range("G169").Value = 0.93910354605886581
$ws.Range("H169").Value = 0.84803759090643216
$ws.Range("I169").Value = 0.15885503698819448
$ws.Range("F173").Value = 0.62124242606729452
$ws.Range("G173").Value = 0.93897764033575126
$ws.Range("H173").Value = 0.84795401402961812
$ws.Range("I173").Value = 0.15881821467293097
$ws.Range("F177").Value = 0.66189198041641251
$ws.Range("G177").Value = 0.93782084188911707
$ws.Range("H177").Value = 0.84828629032258063
$ws.Range("I177").Value = 0.15898963918272713
$ws.Range("F181").Value = 0.68969916772274986
$ws.Range("G181").Value = 0.92110773237446264
$ws.Range("H181").Value = 0.84905088047000976
$ws.Range("I181").Value = 0.16028098990569931
$ws.Range("F185").Value = 0.69349744469145214
$ws.Range("G185").Value = 0.9210610040324072
$ws.Range("H185").Value = 0.84896356934271622
$ws.Range("I185").Value = 0.16022872559975032
$ws.Range("F189").Value = 0.69885968492947825
$ws.Range("G189").Value = 0.92095098561534361
$ws.Range("H189").Value = 0.84894438066079858
$ws.Range("I189").Value = 0.16022077189589881
$ws.Range("F193").Value = 0.7265430694955971
$ws.Range("G193").Value = 0.91983828451263216
$ws.Range("H193").Value = 0.84900076374745415
$ws.Range("I193").Value = 0.16026908181378086
